$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C10").Value = 8200
$ws.Range("C11:C23").Value = 8114
$ws.Range("C24:C31").Value = 8063
$ws.Range("C32:C36").Value = 8019
$ws.Range("C37:C44").Value = 7994
$ws.Range("C45:C48").Value = 7861
$ws.Range("C49:C69").Value = 7859
$ws.Range("C70:C71").Value = 7785
$ws.Range("C72:C84").Value = 7320
$ws.Range("C85:C88").Value = 7295
